# Updated mapping of SubjectID and RecordID
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Subject ID row (row 15) - NIEM 3.1 Mapping column C
$ws.Range("C15").Value = "/cdr-doc:ConsentDecisionReport/j:Booking/j:BookingSubject/j:SubjectIdentification/nc:IdentificationID"

# Booking Number row (row 16) - NIEM 3.1 Mapping column C
$ws.Range("C16").Value = "cdr-doc:ConsentDecisionReport/j:Booking/j:BookingAgencyRecordIdentification/nc:IdentificationID"

# Update the active selection to reflect where editing ended up
$ws.Range("C16").Select()
